$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with matching run formatting) ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Cells whose style must switch between the text placeholder style (s=14)
# and a numeric style (s=15 / s=16); copy the format from a stable anchor cell
# of the desired style first, then set the new value. ---
$ws.Range("L15").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = -100
$ws.Range("L15").Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("M15").Value = -100
$ws.Range("C16").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 5
$ws.Range("A14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("L15").Copy()
$ws.Range("L20").PasteSpecial(-4122)
$ws.Range("L20").Value = 400
$ws.Range("A14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "0"
$ws.Range("C16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("L15").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("C16").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = 1
$ws.Range("L15").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = -100
$ws.Range("C16").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("J28").Value = 1
$ws.Range("L15").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = -100
$ws.Range("C16").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("L15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("C16").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("L15").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = -100
$ws.Range("C16").Copy()
$ws.Range("J29").PasteSpecial(-4122)
$ws.Range("J29").Value = 1
$ws.Range("L15").Copy()
$ws.Range("K29").PasteSpecial(-4122)
$ws.Range("K29").Value = -100
$ws.Range("L15").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = -100

# --- Plain value updates (style/format unchanged) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = 16.666666666666
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -77.419354838709
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 700
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -6.666666666666
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 18.181818181818
$ws.Range("L17").Value = 18.181818181818
$ws.Range("M17").Value = 85.714285714285
$ws.Range("N17").Value = -7.142857142857
$ws.Range("E18").Value = 66.666666666666
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -11.111111111111
$ws.Range("M18").Value = -46.666666666666
$ws.Range("N18").Value = -85.185185185185
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = 11.111111111111
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 33
$ws.Range("K19").Value = 30.303030303030
$ws.Range("L19").Value = 72
$ws.Range("M19").Value = 53.571428571428
$ws.Range("N19").Value = -46.913580246913
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 100
$ws.Range("N20").Value = -81.481481481481
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 52.941176470588
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = 16.455696202531
$ws.Range("I21").Value = 76
$ws.Range("J21").Value = 61
$ws.Range("K21").Value = 24.590163934426
$ws.Range("L21").Value = 55.102040816326
$ws.Range("M21").Value = 31.034482758620
$ws.Range("N21").Value = -63.636363636363
$ws.Range("C22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 5
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 150
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 2
$ws.Range("L23").Value = 100
$ws.Range("M23").Value = -33.333333333333
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 72.222222222222
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 48.529411764705
$ws.Range("I24").Value = 68
$ws.Range("J24").Value = 52
$ws.Range("K24").Value = 30.769230769230
$ws.Range("L24").Value = 13.333333333333
$ws.Range("M24").Value = 4.615384615384
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("I25").Value = 22
$ws.Range("J25").Value = 17
$ws.Range("K25").Value = 29.411764705882
$ws.Range("L25").Value = 100
$ws.Range("M25").Value = 69.230769230769
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 100
$ws.Range("L27").Value = 200
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 2
